# Updates cryptos list values (Price/Volume columns) per the Fri Aug 25 13:23:09 UTC 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.220.13"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.662.44"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'217.36"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Value = "'0.5226"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'0.2645"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "'0.06286"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").Value = "'20.86"
$ws.Range("E10").Value = "  -3.84%  "
$ws.Range("D11").Value = "'0.07766"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.632.11"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("D14").Value = "1.887.40"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "'0.5464"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "0.0₅8173"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("D17").Value = "'64.99"
$ws.Range("D18").Value = "26.222.17"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'4.606"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").Value = "'192.11"
$ws.Range("E21").Value = "  -0.56%  "
$ws.Range("D22").Value = "'10.04"
$ws.Range("E22").Value = "  -2.96%  "
$ws.Range("D23").Value = "'6.017"
$ws.Range("E23").Value = "  -4.03%  "
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").Value = "'139.19"
$ws.Range("E25").Value = "  +0.43%  "
$ws.Range("D26").Value = "'0.1229"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("D27").Value = "'7.281"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").Value = "'16.20"
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("D29").Value = "'1.433"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("D30").Value = "'0.05974"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").Value = "'1.275"
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").Value = "'3.543"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "'3.272"
$ws.Range("E33").Value = "  -3.66%  "
$ws.Range("D34").Value = "'1.583"
$ws.Range("E34").Value = "  -5.86%  "
$ws.Range("D35").Value = "'0.9613"
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("D36").Value = "'2.414"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").Value = "'2.772"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'0.5693"
$ws.Range("E38").Value = "  -5.96%  "
$ws.Range("D39").Value = "'0.01598"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").Value = "'5.975"
$ws.Range("E40").Value = "  -0.87%  "
$ws.Range("D41").Value = "'0.8494"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "1.005.33"
$ws.Range("E43").Value = "  -8.15%  "
$ws.Range("D44").Value = "'100.52"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "1.801.61"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("D46").Value = "0.0₈109"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").Value = "'56.63"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").Value = "'8.050"
$ws.Range("E49").Value = "  -1.59%  "
$ws.Range("D50").Value = "'0.4342"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").Value = "'0.05151"
$ws.Range("E51").Value = "  -1.06%  "
